# Update the "取得日時" (acquisition timestamp) column on the "ランサーズ" sheet
# for all data rows (2-11) to reflect the new run time: 2026-01-20 12:58:58

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-20 12:58:58"

for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
